# Algeria Ligue 1 2023-2024 update script
# - Swaps the match-detail columns (F:V) between several row pairs that
#   were previously out of chronological order within their matchday.
# - Appends 5 new match rows (82-86) for the 2024-01-03 / 2024-01-04 / 2024-01-05 fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Swap F:V between the paired rows.
# ---------------------------------------------------------------------
$pairs = @(
    @(12, 13),
    @(21, 22),
    @(23, 24),
    @(31, 32),
    @(41, 42),
    @(51, 52),
    @(64, 65)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $range1 = $ws.Range("F" + $row1 + ":V" + $row1)
    $range2 = $ws.Range("F" + $row2 + ":V" + $row2)

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}

# ---------------------------------------------------------------------
# 2. Append the 5 new rows (82-86), copying the formatting of the
#    last existing data row (81) and then filling in the values.
# ---------------------------------------------------------------------
$ws.Range("A81:V81").Copy()
$ws.Range("A82:V86").PasteSpecial(-4122)

$newRows = @(
    @{ Row=82; Indice=81; E=45296.63541666666; F="Ben Aknoun"; G=2; H="Oran";      I=1;
       J=2.42; K="04/01/2024 07:10"; L=2.5;  M="05/01/2024 14:56";
       N=2.75; O="04/01/2024 07:10"; P=2.73; Q="05/01/2024 15:06";
       R=3.49; S="04/01/2024 07:10"; T=3.43; U="05/01/2024 14:58";
       V="https://www.betexplorer.com/football/algeria/ligue-1/es-ben-aknoun-oran/EDkdZimb/" },

    @{ Row=83; Indice=82; E=45296.63541666666; F="Magra";      G=0; H="ES Setif";  I=1;
       J=2.15; K="03/01/2024 08:01"; L=2.92; M="05/01/2024 15:01";
       N=3.3;  O="03/01/2024 08:01"; P=2.98; Q="05/01/2024 15:01";
       R=3.09; S="03/01/2024 08:01"; T=2.64; U="05/01/2024 14:39";
       V="https://www.betexplorer.com/football/algeria/ligue-1/magra-es-setif/YaHZHjIN/" },

    @{ Row=84; Indice=83; E=45296.63541666666; F="US Souf";    G=1; H="Paradou";   I=4;
       J=2.92; K="03/01/2024 15:42"; L=4.36; M="05/01/2024 15:10";
       N=2.7;  O="03/01/2024 15:42"; P=3.08; Q="05/01/2024 15:10";
       R=2.63; S="03/01/2024 15:42"; T=1.99; U="05/01/2024 15:10";
       V="https://www.betexplorer.com/football/algeria/ligue-1/us-souf-paradou/f9IVIW2H/" },

    @{ Row=85; Indice=84; E=45296.65625;       F="ASO Chlef";  G=0; H="USM Alger"; I=1;
       J=2.09; K="04/01/2024 07:12"; L=1.85; M="05/01/2024 15:20";
       N=2.85; O="04/01/2024 07:12"; P=3.18; Q="05/01/2024 15:40";
       R=3.8;  S="04/01/2024 07:12"; T=4.9;  U="05/01/2024 15:20";
       V="https://www.betexplorer.com/football/algeria/ligue-1/aso-chlef-usm-alger/SMGwHAXT/" },

    @{ Row=86; Indice=85; E=45296.77083333334; F="Constantine"; G=1; H="Biskra";   I=1;
       J=1.48; K="03/01/2024 18:42"; L=1.35; M="05/01/2024 18:25";
       N=3.7;  O="03/01/2024 18:42"; P=4.53; Q="05/01/2024 18:27";
       R=7.01; S="03/01/2024 18:42"; T=10.68; U="05/01/2024 18:27";
       V="https://www.betexplorer.com/football/algeria/ligue-1/constantine-biskra/pvohzYXi/" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Indice
    $ws.Cells.Item($row, 2).Value = "algeria"
    $ws.Cells.Item($row, 3).Value = "ligue-1"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
}
